$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 80
$ws.Range("H80").Value = 849
$ws.Range("I80").Value = 596.7
$ws.Range("J80").Value = 1078.3636
$ws.Range("K80").Value = 1790.1
$ws.Range("L80").Value = 3235.0908
$ws.Range("M80").Value = -792.1000000000001
$ws.Range("N80").Value = -5231.0908

# Row 83
$ws.Range("H83").Value = 849
$ws.Range("I83").Value = 596.7
$ws.Range("J83").Value = 1078.3636
$ws.Range("K83").Value = 5370.3
$ws.Range("L83").Value = 9705.2724
$ws.Range("M83").Value = -378.3000000000002
$ws.Range("N83").Value = -19689.2724

# Row 113
$ws.Range("H113").Value = 4704.375
$ws.Range("I113").Value = 4028.5
$ws.Range("J113").Value = 5830.8335
$ws.Range("K113").Value = 4028.5
$ws.Range("L113").Value = 5830.8335
$ws.Range("M113").Value = -774.5
$ws.Range("N113").Value = -12338.8335

# Row 132
$ws.Range("H132").Value = 879727.6
$ws.Range("I132").Value = 879727.6
$ws.Range("K132").Value = 2639182.8
$ws.Range("M132").Value = -2636652.8

# Row 137
$ws.Range("H137").Value = 1684.0769
$ws.Range("I137").Value = 1517.6364
$ws.Range("K137").Value = 4552.9092
$ws.Range("M137").Value = -2002.9092

$ws = $wb.Worksheets.Item("ARM")
# Row 2
$ws.Range("H2").Value = 1660.9524
$ws.Range("J2").Value = 1318.4
$ws.Range("L2").Value = 1318.4
$ws.Range("N2").Value = -1544.4

# Row 61
$ws.Range("H61").Value = 4249.2383
$ws.Range("I61").Value = 3032.6
$ws.Range("K61").Value = 3032.6
$ws.Range("M61").Value = -2820.6

# Row 63
$ws.Range("H63").Value = 2134.8333
$ws.Range("J63").Value = 4000
$ws.Range("L63").Value = 4000
$ws.Range("N63").Value = -5372

# Row 66
$ws.Range("H66").Value = 2134.8333
$ws.Range("J66").Value = 4000
$ws.Range("L66").Value = 20000
$ws.Range("N66").Value = -26864

# Row 74
$ws.Range("H74").Value = 1858.7551
$ws.Range("I74").Value = 1715.5526
$ws.Range("J74").Value = 2353.4546
$ws.Range("K74").Value = 1715.5526
$ws.Range("L74").Value = 2353.4546
$ws.Range("M74").Value = -841.5526
$ws.Range("N74").Value = -4101.4546

# Row 77
$ws.Range("H77").Value = 1858.7551
$ws.Range("I77").Value = 1715.5526
$ws.Range("J77").Value = 2353.4546
$ws.Range("K77").Value = 8577.762999999999
$ws.Range("L77").Value = 11767.273
$ws.Range("M77").Value = -4209.762999999999
$ws.Range("N77").Value = -20503.273

# Row 88
$ws.Range("H88").Value = 1649.6364
$ws.Range("I88").Value = 1997.25
$ws.Range("J88").Value = 1451
$ws.Range("K88").Value = 1997.25
$ws.Range("L88").Value = 1451
$ws.Range("M88").Value = -1591.25
$ws.Range("N88").Value = -2263

# Row 91
$ws.Range("H91").Value = 1649.6364
$ws.Range("I91").Value = 1997.25
$ws.Range("J91").Value = 1451
$ws.Range("K91").Value = 1997.25
$ws.Range("L91").Value = 1451
$ws.Range("M91").Value = -593.25
$ws.Range("N91").Value = -4259

# Row 116
$ws.Range("H116").Value = 1660.9524
$ws.Range("J116").Value = 1318.4
$ws.Range("L116").Value = 1318.4
$ws.Range("N116").Value = -5906.4

# Row 136
$ws.Range("H136").Value = 4249.2383
$ws.Range("I136").Value = 3032.6
$ws.Range("K136").Value = 9097.799999999999
$ws.Range("M136").Value = -6547.799999999999

$ws = $wb.Worksheets.Item("BSM")
# Row 3
$ws.Range("H3").Value = 1660.9524
$ws.Range("J3").Value = 1318.4
$ws.Range("L3").Value = 1318.4
$ws.Range("N3").Value = -1546.4

# Row 86
$ws.Range("H86").Value = 47619948
$ws.Range("I86").Value = 100000870
$ws.Range("J86").Value = 926.8182
$ws.Range("K86").Value = 100000870
$ws.Range("L86").Value = 926.8182
$ws.Range("M86").Value = -99999747
$ws.Range("N86").Value = -3172.8182

# Row 89
$ws.Range("H89").Value = 47619948
$ws.Range("I89").Value = 100000870
$ws.Range("J89").Value = 926.8182
$ws.Range("K89").Value = 500004350
$ws.Range("L89").Value = 4634.091
$ws.Range("M89").Value = -499998734
$ws.Range("N89").Value = -15866.091

# Row 99
$ws.Range("H99").Value = 1793.8572
$ws.Range("I99").Value = 1808.1666
$ws.Range("J99").Value = 1708
$ws.Range("K99").Value = 1808.1666
$ws.Range("L99").Value = 1708
$ws.Range("M99").Value = -310.1666
$ws.Range("N99").Value = -4704

$ws = $wb.Worksheets.Item("CRP")
# Row 31
$ws.Range("H31").Value = 1654.92
$ws.Range("I31").Value = 1127.6923
$ws.Range("K31").Value = 1127.6923
$ws.Range("M31").Value = -832.6922999999999

# Row 34
$ws.Range("H34").Value = 1654.92
$ws.Range("I34").Value = 1127.6923
$ws.Range("K34").Value = 1127.6923
$ws.Range("M34").Value = -925.6922999999999

# Row 58
$ws.Range("H58").Value = 50003304
$ws.Range("I58").Value = 25002896
$ws.Range("J58").Value = 83337180
$ws.Range("K58").Value = 25002896
$ws.Range("L58").Value = 83337180
$ws.Range("M58").Value = -25002693
$ws.Range("N58").Value = -83337586

# Row 99
$ws.Range("H99").Value = 3390
$ws.Range("I99").Value = 3437
$ws.Range("K99").Value = 3437
$ws.Range("M99").Value = -1939

# Row 126
$ws.Range("H126").Value = 3390
$ws.Range("I126").Value = 3437
$ws.Range("K126").Value = 10311
$ws.Range("M126").Value = -7841

# Row 134
$ws.Range("H134").Value = 13893114
$ws.Range("I134").Value = 18522918
$ws.Range("K134").Value = 55568754
$ws.Range("M134").Value = -55566219

# Row 136
$ws.Range("H136").Value = 50003304
$ws.Range("I136").Value = 25002896
$ws.Range("J136").Value = 83337180
$ws.Range("K136").Value = 75008688
$ws.Range("L136").Value = 250011540
$ws.Range("M136").Value = -75006138
$ws.Range("N136").Value = -250016640

$ws = $wb.Worksheets.Item("CUL")
# Row 76
$ws.Range("H76").Value = 6814
$ws.Range("J76").Value = 6814
$ws.Range("L76").Value = 20442
$ws.Range("N76").Value = -21208

# Row 79
$ws.Range("H79").Value = 6814
$ws.Range("J79").Value = 6814
$ws.Range("L79").Value = 20442
$ws.Range("N79").Value = -23094

# Row 81
$ws.Range("H81").Value = 2334.8333
$ws.Range("I81").Value = 1256.5
$ws.Range("K81").Value = 3769.5
$ws.Range("M81").Value = -2646.5

# Row 82
$ws.Range("H82").Value = 4844
$ws.Range("I82").Value = 971
$ws.Range("K82").Value = 2913
$ws.Range("M82").Value = -2507

# Row 84
$ws.Range("H84").Value = 2334.8333
$ws.Range("I84").Value = 1256.5
$ws.Range("K84").Value = 11308.5
$ws.Range("M84").Value = -5692.5

# Row 85
$ws.Range("H85").Value = 4844
$ws.Range("I85").Value = 971
$ws.Range("K85").Value = 2913
$ws.Range("M85").Value = -1509

# Row 95
$ws.Range("H95").Value = 14527
$ws.Range("J95").Value = 14527
$ws.Range("L95").Value = 43581
$ws.Range("N95").Value = -47699

# Row 96
$ws.Range("H96").Value = 12925.25
$ws.Range("J96").Value = 12925.25
$ws.Range("L96").Value = 38775.75
$ws.Range("N96").Value = -42893.75

# Row 98
$ws.Range("H98").Value = 252
$ws.Range("I98").Value = 0
$ws.Range("J98").Value = 252
$ws.Range("K98").Value = 0
$ws.Range("L98").Value = 756
$ws.Range("M98").ClearContents() | Out-Null
$ws.Range("N98").Value = -3752

# Row 99
$ws.Range("H99").Value = 4331.6665
$ws.Range("J99").Value = 0
$ws.Range("L99").Value = 0
$ws.Range("N99").ClearContents() | Out-Null

# Row 100
$ws.Range("H100").Value = 6664.6665
$ws.Range("I100").Value = 6664.6665
$ws.Range("J100").Value = 0
$ws.Range("K100").Value = 19993.9995
$ws.Range("L100").Value = 0
$ws.Range("M100").Value = -19182.9995
$ws.Range("N100").ClearContents() | Out-Null

# Row 102
$ws.Range("H102").Value = 2999.5
$ws.Range("I102").Value = 2999.5
$ws.Range("J102").Value = 0
$ws.Range("K102").Value = 8998.5
$ws.Range("L102").Value = 0
$ws.Range("M102").Value = -6564.5
$ws.Range("N102").ClearContents() | Out-Null

# Row 104
$ws.Range("H104").Value = 73043.89999999999
$ws.Range("I104").Value = 81111
$ws.Range("J104").Value = 440
$ws.Range("K104").Value = 243333
$ws.Range("L104").Value = 1320
$ws.Range("M104").Value = -240712
$ws.Range("N104").Value = -6562

# Row 139
$ws.Range("H139").Value = 1644.4
$ws.Range("I139").Value = 1000.2143
$ws.Range("K139").Value = 3000.6429
$ws.Range("M139").Value = 2139.3571

$ws = $wb.Worksheets.Item("GSM")
# Row 31
$ws.Range("H31").Value = 7941.857
$ws.Range("I31").Value = 5932.1665
$ws.Range("J31").Value = 20000
$ws.Range("K31").Value = 5932.1665
$ws.Range("L31").Value = 20000
$ws.Range("M31").Value = -5640.1665
$ws.Range("N31").Value = -20584

# Row 37
$ws.Range("H37").Value = 7941.857
$ws.Range("I37").Value = 5932.1665
$ws.Range("J37").Value = 20000
$ws.Range("K37").Value = 5932.1665
$ws.Range("L37").Value = 20000
$ws.Range("M37").Value = -5655.1665
$ws.Range("N37").Value = -20554

# Row 113
$ws.Range("H113").Value = 36670.332
$ws.Range("I113").Value = 100011
$ws.Range("K113").Value = 100011
$ws.Range("M113").Value = -97841

# Row 126
$ws.Range("I126").Value = 76925960
$ws.Range("J126").Value = 7120.75
$ws.Range("K126").Value = 230777880
$ws.Range("L126").Value = 21362.25
$ws.Range("M126").Value = -230775410
$ws.Range("N126").Value = -26302.25

# Row 140
$ws.Range("H140").Value = 49531.25
$ws.Range("I140").Value = 49531.25
$ws.Range("J140").Value = 0
$ws.Range("K140").Value = 49531.25
$ws.Range("L140").Value = 0
$ws.Range("M140").Value = -44351.25
$ws.Range("N140").ClearContents() | Out-Null

$ws = $wb.Worksheets.Item("LTW")
# Row 40
$ws.Range("H40").Value = 23940.75
$ws.Range("J40").Value = 4632.6665
$ws.Range("L40").Value = 4632.6665
$ws.Range("N40").Value = -4904.6665

# Row 128
$ws.Range("H128").Value = 69394.5
$ws.Range("J128").Value = 69394.5
$ws.Range("L128").Value = 69394.5
$ws.Range("N128").Value = -79354.5

$ws = $wb.Worksheets.Item("WVR")
# Row 113
$ws.Range("H113").Value = 1819.579
$ws.Range("I113").Value = 1799.1177
$ws.Range("K113").Value = 5397.3531
$ws.Range("M113").Value = -3227.3531

# Row 122
$ws.Range("H122").Value = 3277.077
$ws.Range("I122").Value = 2843.5908
$ws.Range("K122").Value = 8530.7724
$ws.Range("M122").Value = -6080.7724
